$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (row 2), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Update the view: zoom level and selected cell.
$ws.Activate()
$excel.ActiveWindow.Zoom = 92
$ws.Range("I16").Select()
